# Append new barcode rows (5-8) to sheet1, matching existing column formatting:
#   Column A: barcode (integer number format, inherited from column style)
#   Column B: quantity (general/default format)
#   Column C: expiry date (date number format, same as existing rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(3502110008091, 20, 44776),
    @(5410013110002, 5, 44596),
    @(5411028070480, 90, 44590),
    @(5411188115472, 63, 44589)
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]

    # Copy format from the cell directly above (reuses the existing date
    # style instead of minting a brand-new numFmt entry).
    $ws.Cells.Item($r - 1, 3).Copy() | Out-Null
    $dateCell = $ws.Cells.Item($r, 3)
    $dateCell.PasteSpecial(-4122) | Out-Null
    $dateCell.Value = $data[2]
}

# Update selection to match post-edit state (cursor moved to next empty row)
$ws.Range("A9").Select() | Out-Null
